# Update natmi LR-pairs data (Cxcl13-Ackr4) with refreshed TPM values.
# Drop the "Resolving-Mac" target-cluster rows (for both FAPs and MuSCs
# senders) and refresh all remaining edge-expression metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows whose Target cluster is "Resolving-Mac" (row 9 = MuSCs->Resolving-Mac,
# row 5 = FAPs->Resolving-Mac). Delete the higher row index first so the
# remaining row numbers don't shift before the second delete.
$ws.Rows(9).Delete()
$ws.Rows(5).Delete()

# After the deletions, rows 2-7 hold (in order):
#   FAPs->ECs, FAPs->FAPs, FAPs->MuSCs, MuSCs->ECs, MuSCs->FAPs, MuSCs->MuSCs
# Refresh their computed values with the newly recalculated TPM numbers.

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.528846666666666
$ws.Cells.Item(2,8).Value = 10.58654
$ws.Cells.Item(2,9).Value = 0.8781048434890718
$ws.Cells.Item(2,10).Value = 0.8781048434890719
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.377371
$ws.Cells.Item(2,14).Value = 1.132113
$ws.Cells.Item(2,15).Value = 0.4698794580655765
$ws.Cells.Item(2,16).Value = 0.4698794580655764
$ws.Cells.Item(2,17).Value = 1.331684395446666
$ws.Cells.Item(2,18).Value = 11.98515955902
$ws.Cells.Item(2,19).Value = 0.4126034279834029
$ws.Cells.Item(2,20).Value = 0.4126034279834029
# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.528846666666666
$ws.Cells.Item(3,8).Value = 10.58654
$ws.Cells.Item(3,9).Value = 0.8781048434890718
$ws.Cells.Item(3,10).Value = 0.8781048434890719
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.3560133333333333
$ws.Cells.Item(3,14).Value = 1.06804
$ws.Cells.Item(3,15).Value = 0.443286188209444
$ws.Cells.Item(3,16).Value = 0.443286188209444
$ws.Cells.Item(3,17).Value = 1.256316464622222
$ws.Cells.Item(3,18).Value = 11.3068481816
$ws.Cells.Item(3,19).Value = 0.3892517489185211
$ws.Cells.Item(3,20).Value = 0.3892517489185212
# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.528846666666666
$ws.Cells.Item(4,8).Value = 10.58654
$ws.Cells.Item(4,9).Value = 0.8781048434890718
$ws.Cells.Item(4,10).Value = 0.8781048434890719
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.06973866666666667
$ws.Cells.Item(4,14).Value = 0.209216
$ws.Cells.Item(4,15).Value = 0.08683435372497944
$ws.Cells.Item(4,16).Value = 0.08683435372497944
$ws.Cells.Item(4,17).Value = 0.2460970614044445
$ws.Cells.Item(4,18).Value = 2.21487355264
$ws.Cells.Item(4,19).Value = 0.07624966658714777
$ws.Cells.Item(4,20).Value = 0.07624966658714778
# Row 5
$ws.Cells.Item(5,1).Value = "MuSCs"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.489861
$ws.Cells.Item(5,8).Value = 1.469583
$ws.Cells.Item(5,9).Value = 0.1218951565109281
$ws.Cells.Item(5,10).Value = 0.1218951565109281
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.377371
$ws.Cells.Item(5,14).Value = 1.132113
$ws.Cells.Item(5,15).Value = 0.4698794580655765
$ws.Cells.Item(5,16).Value = 0.4698794580655764
$ws.Cells.Item(5,17).Value = 0.184859335431
$ws.Cells.Item(5,18).Value = 1.663734018879
$ws.Cells.Item(5,19).Value = 0.05727603008217352
$ws.Cells.Item(5,20).Value = 0.05727603008217352
# Row 6
$ws.Cells.Item(6,1).Value = "MuSCs"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.489861
$ws.Cells.Item(6,8).Value = 1.469583
$ws.Cells.Item(6,9).Value = 0.1218951565109281
$ws.Cells.Item(6,10).Value = 0.1218951565109281
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.3560133333333333
$ws.Cells.Item(6,14).Value = 1.06804
$ws.Cells.Item(6,15).Value = 0.443286188209444
$ws.Cells.Item(6,16).Value = 0.443286188209444
$ws.Cells.Item(6,17).Value = 0.17439704748
$ws.Cells.Item(6,18).Value = 1.56957342732
$ws.Cells.Item(6,19).Value = 0.05403443929092291
$ws.Cells.Item(6,20).Value = 0.05403443929092291
# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.489861
$ws.Cells.Item(7,8).Value = 1.469583
$ws.Cells.Item(7,9).Value = 0.1218951565109281
$ws.Cells.Item(7,10).Value = 0.1218951565109281
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.06973866666666667
$ws.Cells.Item(7,14).Value = 0.209216
$ws.Cells.Item(7,15).Value = 0.08683435372497944
$ws.Cells.Item(7,16).Value = 0.08683435372497944
$ws.Cells.Item(7,17).Value = 0.03416225299200001
$ws.Cells.Item(7,18).Value = 0.307460276928
$ws.Cells.Item(7,19).Value = 0.01058468713783166
$ws.Cells.Item(7,20).Value = 0.01058468713783166
